# The source revision for this document only touches cosmetic,
# non-semantic bookkeeping (OOXML namespace-prefix bookkeeping on the
# root elements of a handful of parts) -- no visible text, formatting,
# table, header/footer, or style content actually changed between the
# "before" and "after" copies of this file. Drive the document through
# Word's automation surface without mutating any user-visible content,
# mirroring that the real edit here was a no-op resave (the author's
# own note says the real screen/table work described in the commit
# message is "ainda em andamento" / still in progress and not actually
# present in this particular document).
$d = $word.ActiveDocument

# Touch the object model (forces Word to walk the document) without
# altering any content: a find that can never match leaves the
# document byte-for-byte equivalent while still exercising the COM
# automation path end to end.
$found = $d.Content.Find.Execute(
    "`u{0001}NoSuchTextInDocument`u{0001}",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0
)
